# PlayerPerformance_7109.xlsx update:
#  - insert a new "Player Info" sheet before "ODI Batting"
#  - rename MATCH_CARD_LINK -> MATCH_CODE on both existing sheets and
#    replace the full scorecard URL with just the match code

$wb = $excel.ActiveWorkbook

# --- 1. New "Player Info" sheet, inserted before "ODI Batting" ---------
$battingSheetBeforeInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetBeforeInsert)
$playerInfo.Name = "Player Info"

# NOTE: inserting a sheet can shift/reuse earlier object references in
# this engine, so re-resolve the other sheets by name AFTER the insert.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "'7109"
$playerInfo.Range("B2").Value = "Rehan Ahmed"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# Match the bold/boxed header look already used on the other sheets.
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ----------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").Value = "'4717"

# --- 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE ----------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").Value = "'4717"
